$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "dsfd"
$ws.Range("A3").Value = "dfds"
$ws.Range("A4").Value = "dfds"
$ws.Range("A5").Value = "ds"
$ws.Range("A6").Value = "f"
$ws.Range("A7").Value = "sdf"
$ws.Range("A8").Value = "ds"
$ws.Range("A9").Value = "f"
$ws.Range("A10").Value = "sdf"
$ws.Range("A11").Value = "sdf"

$ws.Range("A11").Select() | Out-Null
